{"js": "// Word JS API (Office.js) edit script.\n// Body of: async (context) => { ... }\n//\n// Change: the sentence \"... protected and can't be delated. \" loses its\n// trailing period and gains a parenthetical explanation before a new\n// closing period, i.e.:\n//   \" protected and can't be delated. \"\n// becomes\n//   \" protected and can't be delated (It is hardcoded to database and\n//     there is no setter metod for this property). \"\n\nconst body = context.document.body;\n\nconst oldText = \" protected and can't be delated. \";\nconst newText =\n  \" protected and can't be delated\" +\n  \" (It is hardcoded to database and there is no setter metod for this property)\" +\n  \". \";\n\nconst searchResults = body.search(oldText, { matchCase: true, matchWholeWord: false });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Target sentence fragment not found: \" + oldText);\n}\n\n// Replace the exact matched range in place so the surrounding runs\n// (and their formatting) are left untouched.\nsearchResults.items[0].insertText(newText, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word COM interop edit script.\n# $word.ActiveDocument is the open document.\n#\n# Change: the sentence \"... protected and can't be delated. \" loses its\n# trailing period and gains a parenthetical explanation before a new\n# closing period, i.e.:\n#   \" protected and can't be delated. \"\n# becomes\n#   \" protected and can't be delated (It is hardcoded to database and\n#     there is no setter metod for this property). \"\n\n$d = $word.ActiveDocument\n\n$oldText = \" protected and can't be delated. \"\n$newText = \" protected and can't be delated (It is hardcoded to database and there is no setter metod for this property). \"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = $oldText\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.Wrap = 0  # wdFindStop\n$found = $find.Execute()\n\nif ($found -and $find.Found) {\n    # $find.Parent is the exact matched Range; replacing its .Text in\n    # place edits only that span of the paragraph.\n    $rng = $find.Parent\n    $rng.Text = $newText\n} else {\n    throw \"Target sentence fragment not found: $oldText\"\n}\n"}
